$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.04191
$ws.Range("H2").Value = 0.12573
$ws.Range("I2").Value = 0.0003002276973850376
$ws.Range("J2").Value = 0.0003002276973850376
$ws.Range("M2").Value = 7.487621999999999
$ws.Range("N2").Value = 22.462866
$ws.Range("O2").Value = 0.1384395179233961
$ws.Range("P2").Value = 0.1384395179233961
$ws.Range("Q2").Value = 0.31380623802
$ws.Range("R2").Value = 2.82425614218
$ws.Range("S2").Value = 0.00004156337769323584
$ws.Range("T2").Value = 0.00004156337769323584

$ws.Range("G3").Value = 0.04191
$ws.Range("H3").Value = 0.12573
$ws.Range("I3").Value = 0.0003002276973850376
$ws.Range("J3").Value = 0.0003002276973850376
$ws.Range("O3").Value = 0.5916411627275552
$ws.Range("P3").Value = 0.5916411627275552
$ws.Range("Q3").Value = 1.34109602748
$ws.Range("R3").Value = 12.06986424732
$ws.Range("S3").Value = 0.0001776270639639002
$ws.Range("T3").Value = 0.0001776270639639002

$ws.Range("G4").Value = 0.04191
$ws.Range("H4").Value = 0.12573
$ws.Range("I4").Value = 0.0003002276973850376
$ws.Range("J4").Value = 0.0003002276973850376
$ws.Range("M4").Value = 14.59882166666667
$ws.Range("N4").Value = 43.796465
$ws.Range("O4").Value = 0.2699193193490487
$ws.Range("P4").Value = 0.2699193193490487
$ws.Range("Q4").Value = 0.6118366160500001
$ws.Range("R4").Value = 5.50652954445
$ws.Range("S4").Value = 0.00008103725572790152
$ws.Range("T4").Value = 0.00008103725572790152

$ws.Range("I5").Value = 0.9976864582107258
$ws.Range("J5").Value = 0.9976864582107259
$ws.Range("M5").Value = 7.487621999999999
$ws.Range("N5").Value = 22.462866
$ws.Range("O5").Value = 0.1384395179233961
$ws.Range("P5").Value = 0.1384395179233961
$ws.Range("Q5").Value = 1042.809297414972
$ws.Range("R5").Value = 9385.283676734747
$ws.Range("S5").Value = 0.1381192323133933
$ws.Range("T5").Value = 0.1381192323133933

$ws.Range("I6").Value = 0.9976864582107258
$ws.Range("J6").Value = 0.9976864582107259
$ws.Range("O6").Value = 0.5916411627275552
$ws.Range("P6").Value = 0.5916411627275552
$ws.Range("S6").Value = 0.5902723761733302
$ws.Range("T6").Value = 0.5902723761733303

$ws.Range("I7").Value = 0.9976864582107258
$ws.Range("J7").Value = 0.9976864582107259
$ws.Range("M7").Value = 14.59882166666667
$ws.Range("N7").Value = 43.796465
$ws.Range("O7").Value = 0.2699193193490487
$ws.Range("P7").Value = 0.2699193193490487
$ws.Range("Q7").Value = 2033.193845162474
$ws.Range("R7").Value = 18298.74460646227
$ws.Range("S7").Value = 0.2692948497240022
$ws.Range("T7").Value = 0.2692948497240022

$ws.Range("G8").Value = 0.2810466666666667
$ws.Range("H8").Value = 0.84314
$ws.Range("I8").Value = 0.002013314091889132
$ws.Range("J8").Value = 0.002013314091889132
$ws.Range("M8").Value = 7.487621999999999
$ws.Range("N8").Value = 22.462866
$ws.Range("O8").Value = 0.1384395179233961
$ws.Range("P8").Value = 0.1384395179233961
$ws.Range("Q8").Value = 2.10437120436
$ws.Range("R8").Value = 18.93934083924
$ws.Range("S8").Value = 0.0002787222323095114
$ws.Range("T8").Value = 0.0002787222323095114

$ws.Range("G9").Value = 0.2810466666666667
$ws.Range("H9").Value = 0.84314
$ws.Range("I9").Value = 0.002013314091889132
$ws.Range("J9").Value = 0.002013314091889132
$ws.Range("O9").Value = 0.5916411627275552
$ws.Range("P9").Value = 0.5916411627275552
$ws.Range("Q9").Value = 8.99333257464
$ws.Range("R9").Value = 80.93999317175999
$ws.Range("S9").Value = 0.001191159490261058
$ws.Range("T9").Value = 0.001191159490261058

$ws.Range("G10").Value = 0.2810466666666667
$ws.Range("H10").Value = 0.84314
$ws.Range("I10").Value = 0.002013314091889132
$ws.Range("J10").Value = 0.002013314091889132
$ws.Range("M10").Value = 14.59882166666667
$ws.Range("N10").Value = 43.796465
$ws.Range("O10").Value = 0.2699193193490487
$ws.Range("P10").Value = 0.2699193193490487
$ws.Range("Q10").Value = 4.102950166677777
$ws.Range("R10").Value = 36.9265515001
$ws.Range("S10").Value = 0.0005434323693185627
$ws.Range("T10").Value = 0.0005434323693185627
